$wb = $excel.ActiveWorkbook

# ---- Metadata sheet updates ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.4.0-snapshot-1"
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---- Elements sheet updates ----
# The two mapping columns (AK = "Mapping: RIM Mapping", AL = "Mapping:
# Spécification métier vers l'extension ROR AvailableTimeNumberDaysofWeek")
# were swapped: header text and cell values exchange places between the
# two columns, and their (bestFit) column widths follow the content.
$elem = $wb.Worksheets.Item("Elements")

$akHeader = $elem.Range("AK1").Value2
$alHeader = $elem.Range("AL1").Value2
$elem.Range("AK1").Value2 = $alHeader
$elem.Range("AL1").Value2 = $akHeader

for ($r = 2; $r -le 6; $r++) {
    $akCell = $elem.Cells.Item($r, 37)
    $alCell = $elem.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# Swap the bestFit column widths to match the swapped content (column 37
# now holds the long French mapping text, column 38 the short one).
$elem.Columns.Item(37).ColumnWidth = 89.83
$elem.Columns.Item(38).ColumnWidth = 24.15
